$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix sheet/tab name typo: "mars" -> "March"
$ws.Name = "March"

# Ensure Price column values are stored as real numbers instead of text
$ws.Range("D2").Value = 450.0
$ws.Range("D3").Value = 500.0
